# New query for LAST_POKEMON_BY_TYPE
# - Check the "Find the last Pokemon (by National ID) of each type" box
#   (row 13) by setting its linked cell A13 to TRUE, which also flips the
#   associated Check Box 7 form control to its checked state.
# - Move the active selection to B4 (from C6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Tick the checkbox linked to $A$13 ("Find the last Pokemon (by National
# ID) of each type") by writing TRUE to its linked cell.
$ws.Range("A13").Value = $true

# Update the active cell / selection shown when the sheet is opened.
$ws.Range("B4").Select()
